$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("design")

# Row 8 (item #2) corresponds to designator "C11 C12 C21 C22 C31 C32 C41 C42":
# correct the Mfg Part #, Description/Value and Package/Footprint to the 1206 part,
# and clear the note that explained the mismatch.
$ws.Range("G8").Value = "1206"
$ws.Range("F8").Value = "CAP CER 1UF 25V X7R 1206"
$ws.Range("E8").Value = "GCM31CL81H105KA55L"
$ws.Range("K8").ClearContents()

# row height auto-fits to the now-shorter content
$ws.Rows.Item(8).RowHeight = 45.75

# update the active selection to match the authored file
$ws.Range("L9").Select()
